# Updated translation evaluation results: refresh Filename (D) and Impact (E)
# columns for rows 2-82 on the active worksheet to reflect the latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 'codeforces_99_A.py'
$ws.Cells.Item(2, 5).Value = 'Compilation Error'
$ws.Cells.Item(3, 4).Value = 'atcoder_AGC007_C.py'
$ws.Cells.Item(3, 5).Value = 'Compilation Error'
$ws.Cells.Item(4, 4).Value = 'codeforces_373_B.py'
$ws.Cells.Item(4, 5).Value = 'Compilation Error'
$ws.Cells.Item(5, 4).Value = 'atcoder_ABC051_A.py'
$ws.Cells.Item(5, 5).Value = 'Runtime Error'
$ws.Cells.Item(6, 4).Value = 'codeforces_189_A.py'
$ws.Cells.Item(6, 5).Value = 'Runtime Error'
$ws.Cells.Item(7, 4).Value = 'atcoder_ABC122_D.py'
$ws.Cells.Item(7, 5).Value = 'Runtime Error'
$ws.Cells.Item(8, 4).Value = 'codeforces_203_A.py'
$ws.Cells.Item(8, 5).Value = 'Runtime Error'
$ws.Cells.Item(9, 4).Value = 'atcoder_ABC169_C.py'
$ws.Cells.Item(9, 5).Value = 'Runtime Error'
$ws.Cells.Item(10, 4).Value = 'codeforces_669_A.py'
$ws.Cells.Item(10, 5).Value = 'Runtime Error'
$ws.Cells.Item(11, 4).Value = 'atcoder_ABC170_A.py'
$ws.Cells.Item(11, 5).Value = 'Runtime Error'
$ws.Cells.Item(12, 4).Value = 'codeforces_569_A.py'
$ws.Cells.Item(12, 5).Value = 'Test Failed'
$ws.Cells.Item(13, 4).Value = 'atcoder_ABC169_D.py'
$ws.Cells.Item(13, 5).Value = 'Test Failed'
$ws.Cells.Item(14, 4).Value = 'codeforces_59_A.py'
$ws.Cells.Item(14, 5).Value = 'Test Failed'
$ws.Cells.Item(15, 4).Value = 'codeforces_92_A.py'
$ws.Cells.Item(15, 5).Value = 'Test Failed'
$ws.Cells.Item(16, 4).Value = 'codeforces_622_A.py'
$ws.Cells.Item(16, 5).Value = 'Test Failed'
$ws.Cells.Item(17, 4).Value = 'atcoder_ABC178_A.py'
$ws.Cells.Item(17, 5).Value = 'Test Failed'
$ws.Cells.Item(18, 4).Value = 'atcoder_ABC143_A.py'
$ws.Cells.Item(18, 5).Value = 'Test Failed'
$ws.Cells.Item(19, 4).Value = 'codeforces_334_A.py'
$ws.Cells.Item(19, 5).Value = 'Test Failed'
$ws.Cells.Item(20, 4).Value = 'codeforces_514_A.py'
$ws.Cells.Item(20, 5).Value = 'Test Failed'
$ws.Cells.Item(21, 4).Value = 'atcoder_ABC174_C.py'
$ws.Cells.Item(21, 5).Value = 'Test Failed'
$ws.Cells.Item(22, 4).Value = 'codeforces_546_A.py'
$ws.Cells.Item(22, 5).Value = 'Test Failed'
$ws.Cells.Item(23, 4).Value = 'codeforces_110_B.py'
$ws.Cells.Item(23, 5).Value = 'Test Failed'
$ws.Cells.Item(24, 4).Value = 'atcoder_ABC155_E.py'
$ws.Cells.Item(24, 5).Value = 'Test Failed'
$ws.Cells.Item(25, 4).Value = 'codeforces_79_A.py'
$ws.Cells.Item(25, 5).Value = 'Test Failed'
$ws.Cells.Item(26, 4).Value = 'atcoder_AGC046_A.py'
$ws.Cells.Item(26, 5).Value = 'Test Failed'
$ws.Cells.Item(27, 4).Value = 'codeforces_340_A.py'
$ws.Cells.Item(27, 5).Value = 'Test Failed'
$ws.Cells.Item(28, 4).Value = 'atcoder_ARC062_B.py'
$ws.Cells.Item(28, 5).Value = 'Test Failed'
$ws.Cells.Item(29, 4).Value = 'atcoder_ABC108_B.py'
$ws.Cells.Item(29, 5).Value = 'Test Failed'
$ws.Cells.Item(30, 4).Value = 'atcoder_AGC046_B.py'
$ws.Cells.Item(30, 5).Value = 'Test Failed'
$ws.Cells.Item(31, 4).Value = 'codeforces_579_A.py'
$ws.Cells.Item(31, 5).Value = 'Test Failed'
$ws.Cells.Item(32, 4).Value = 'atcoder_ABC149_C.py'
$ws.Cells.Item(32, 5).Value = 'Test Failed'
$ws.Cells.Item(33, 4).Value = 'codeforces_32_B.py'
$ws.Cells.Item(33, 5).Value = 'Test Failed'
$ws.Cells.Item(34, 4).Value = 'atcoder_ABC158_A.py'
$ws.Cells.Item(34, 5).Value = 'Test Failed'
$ws.Cells.Item(35, 4).Value = 'atcoder_ABC042_A.py'
$ws.Cells.Item(35, 5).Value = 'Test Failed'
$ws.Cells.Item(36, 4).Value = 'atcoder_ABC139_B.py'
$ws.Cells.Item(36, 5).Value = 'Test Failed'
$ws.Cells.Item(37, 4).Value = 'codeforces_678_A.py'
$ws.Cells.Item(37, 5).Value = 'Test Failed'
$ws.Cells.Item(38, 4).Value = 'codeforces_58_B.py'
$ws.Cells.Item(38, 5).Value = 'Test Failed'
$ws.Cells.Item(39, 4).Value = 'codeforces_672_A.py'
$ws.Cells.Item(39, 5).Value = 'Test Failed'
$ws.Cells.Item(40, 4).Value = 'codeforces_86_A.py'
$ws.Cells.Item(40, 5).Value = 'Test Failed'
$ws.Cells.Item(41, 4).Value = 'atcoder_ABC125_A.py'
$ws.Cells.Item(41, 5).Value = 'Test Failed'
$ws.Cells.Item(42, 4).Value = 'codeforces_49_A.py'
$ws.Cells.Item(42, 5).Value = 'Test Failed'
$ws.Cells.Item(43, 4).Value = 'atcoder_ABC132_F.py'
$ws.Cells.Item(43, 5).Value = 'Test Failed'
$ws.Cells.Item(44, 4).Value = 'codeforces_190_A.py'
$ws.Cells.Item(44, 5).Value = 'Test Failed'
$ws.Cells.Item(45, 4).Value = 'atcoder_ABC178_B.py'
$ws.Cells.Item(45, 5).Value = 'Test Failed'
$ws.Cells.Item(46, 4).Value = 'atcoder_ABC158_B.py'
$ws.Cells.Item(46, 5).Value = 'Test Failed'
$ws.Cells.Item(47, 4).Value = 'codeforces_96_B.py'
$ws.Cells.Item(47, 5).Value = 'Test Failed'
$ws.Cells.Item(48, 4).Value = 'atcoder_ABC124_C.py'
$ws.Cells.Item(48, 5).Value = 'Test Failed'
$ws.Cells.Item(49, 4).Value = 'codeforces_55_A.py'
$ws.Cells.Item(49, 5).Value = 'Test Failed'
$ws.Cells.Item(50, 4).Value = 'atcoder_ABC142_A.py'
$ws.Cells.Item(50, 5).Value = 'Test Failed'
$ws.Cells.Item(51, 4).Value = 'atcoder_AGC025_A.py'
$ws.Cells.Item(51, 5).Value = 'Test Failed'
$ws.Cells.Item(52, 4).Value = 'atcoder_ABC168_C.py'
$ws.Cells.Item(52, 5).Value = 'Test Failed'
$ws.Cells.Item(53, 4).Value = 'atcoder_ABC043_B.py'
$ws.Cells.Item(53, 5).Value = 'Test Failed'
$ws.Cells.Item(54, 4).Value = 'atcoder_ABC127_B.py'
$ws.Cells.Item(54, 5).Value = 'Test Failed'
$ws.Cells.Item(55, 4).Value = 'codeforces_171_A.py'
$ws.Cells.Item(55, 5).Value = 'Test Failed'
$ws.Cells.Item(56, 4).Value = 'atcoder_ABC120_C.py'
$ws.Cells.Item(56, 5).Value = 'Test Failed'
$ws.Cells.Item(57, 4).Value = 'atcoder_ABC132_A.py'
$ws.Cells.Item(57, 5).Value = 'Test Failed'
$ws.Cells.Item(58, 4).Value = 'codeforces_678_B.py'
$ws.Cells.Item(58, 5).Value = 'Test Failed'
$ws.Cells.Item(59, 4).Value = 'codeforces_276_B.py'
$ws.Cells.Item(59, 5).Value = 'Test Failed'
$ws.Cells.Item(60, 4).Value = 'atcoder_ABC124_A.py'
$ws.Cells.Item(60, 5).Value = 'Test Failed'
$ws.Cells.Item(61, 4).Value = 'codeforces_651_A.py'
$ws.Cells.Item(61, 5).Value = 'Test Failed'
$ws.Cells.Item(62, 4).Value = 'atcoder_ABC153_A.py'
$ws.Cells.Item(62, 5).Value = 'Test Failed'
$ws.Cells.Item(63, 4).Value = 'atcoder_ABC149_B.py'
$ws.Cells.Item(63, 5).Value = 'Test Failed'
$ws.Cells.Item(64, 4).Value = 'codeforces_306_A.py'
$ws.Cells.Item(64, 5).Value = 'Test Failed'
$ws.Cells.Item(65, 4).Value = 'codeforces_544_B.py'
$ws.Cells.Item(65, 5).Value = 'Test Failed'
$ws.Cells.Item(66, 4).Value = 'codeforces_242_A.py'
$ws.Cells.Item(66, 5).Value = 'Test Failed'
$ws.Cells.Item(67, 4).Value = 'codeforces_459_A.py'
$ws.Cells.Item(67, 5).Value = 'Test Failed'
$ws.Cells.Item(68, 4).Value = 'atcoder_ABC136_B.py'
$ws.Cells.Item(68, 5).Value = 'Test Failed'
$ws.Cells.Item(69, 4).Value = 'atcoder_ABC151_A.py'
$ws.Cells.Item(69, 5).Value = 'Test Failed'
$ws.Cells.Item(70, 4).Value = 'atcoder_ABC070_B.py'
$ws.Cells.Item(70, 5).Value = 'Test Failed'
$ws.Cells.Item(71, 4).Value = 'codeforces_379_A.py'
$ws.Cells.Item(71, 5).Value = 'Test Failed'
$ws.Cells.Item(72, 4).Value = 'codeforces_581_A.py'
$ws.Cells.Item(72, 5).Value = 'Test Failed'
$ws.Cells.Item(73, 4).Value = 'atcoder_ABC164_A.py'
$ws.Cells.Item(73, 5).Value = 'Test Failed'
$ws.Cells.Item(74, 4).Value = 'codeforces_8_B.py'
$ws.Cells.Item(74, 5).Value = 'Test Failed'
$ws.Cells.Item(75, 4).Value = 'codeforces_369_B.py'
$ws.Cells.Item(75, 5).Value = 'Test Failed'
$ws.Cells.Item(76, 4).Value = 'atcoder_ABC114_C.py'
$ws.Cells.Item(76, 5).Value = 'Test Failed'
$ws.Cells.Item(77, 4).Value = 'atcoder_AGC002_A.py'
$ws.Cells.Item(77, 5).Value = 'Test Failed'
$ws.Cells.Item(78, 4).Value = 'codeforces_30_A.py'
$ws.Cells.Item(78, 5).Value = 'Test Failed'
$ws.Cells.Item(79, 4).Value = 'atcoder_ARC102_C.py'
$ws.Cells.Item(79, 5).Value = 'Infinite Loop'
$ws.Cells.Item(80, 4).Value = 'atcoder_ABC172_D.py'
$ws.Cells.Item(80, 5).Value = 'Infinite Loop'
$ws.Cells.Item(81, 4).Value = 'atcoder_AGC006_B.py'
$ws.Cells.Item(81, 5).Value = 'Correct'
$ws.Cells.Item(82, 4).Value = 'codeforces_147_A.py'
$ws.Cells.Item(82, 5).Value = 'Correct'
